# jury excel update
# - Rename the single worksheet from "Diretores de departamento" to "jury"
# - Apply the underline cell style (already used by C3/C8) to C2 as well

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "jury"
$ws.Range("C2").Font.Underline = $true
